# Commit: "adding gross_annual_revenue_flag and gross_annual_revenue fields
# to example_sblar.xlsx file"
#
# Adds two new trailing columns (AA: gross_annual_revenue_flag,
# AB: gross_annual_revenue) to the "invalid" sheet, which is already the
# active sheet in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("AA1").Value = "gross_annual_revenue_flag"
$ws.Range("AB1").Value = "gross_annual_revenue"

# --- Data rows (rows 2-11) ---
$ws.Range("AA2").Value = 900
$ws.Range("AB2").Value = 100000

$ws.Range("AA3").Value = 988

$ws.Range("AA4").Value = 900

$ws.Range("AA5").Value = 988
$ws.Range("AB5").Value = 300000

$ws.Range("AB6").Value = 200000

$ws.Range("AA7").Value = 999

$ws.Range("AA8").Value = 990
$ws.Range("AB8").Value = 50000

$ws.Range("AA9").Value = 900
$ws.Range("AB9").Value = 45000

$ws.Range("AA10").Value = 988

$ws.Range("AA11").Value = 900
$ws.Range("AB11").Value = 50000

# Header row grows from two wrapped lines (ht=34) to three (ht=51) now that
# a longer label has been added alongside the existing headers.
$ws.Rows.Item(1).RowHeight = 51

# Reflect the new selection/active cell over the freshly added columns.
$ws.Range("AA1:AB11").Select()
